# Updates cryptos list values (price + 1h volume change) per the commit diff.
# Rows 28/29 (NEARProtocol / RenderToken) are also swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.076.22'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '3.320.28'
$ws.Range("E3").Value = '  +6.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.32'
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.49'
$ws.Range("E6").Value = '  +6.79%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.313.65'
$ws.Range("E8").Value = '  +6.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("E10").Value = '  +3.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.47'
$ws.Range("E11").Value = '  +2.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.475'
$ws.Range("E12").Value = '  +4.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.19'
$ws.Range("E14").Value = '  +3.91%  '
$ws.Range("D15").Value = '3.835.11'
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").Value = '3.294.40'
$ws.Range("E17").Value = '  +5.52%  '
$ws.Range("D18").Value = '64.128.41'
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.95'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '485.43'
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.41'
$ws.Range("E21").Value = '  +2.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.749'
$ws.Range("E22").Value = '  +7.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.10'
$ws.Range("E23").Value = '  +6.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.67'
$ws.Range("E24").Value = '  +5.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.08'
$ws.Range("E25").Value = '  -2.35%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  +3.06%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.33'
$ws.Range("E28").Value = '  +4.74%  '
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.32'
$ws.Range("E29").Value = '  +2.50%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.19'
$ws.Range("E31").Value = '  +7.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.41'
$ws.Range("E32").Value = '  +4.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.109'
$ws.Range("E33").Value = '  +1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.59'
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("E35").Value = '  +2.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.06'
$ws.Range("E36").Value = '  +4.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.17'
$ws.Range("E37").Value = '  +2.21%  '
$ws.Range("D38").Value = '0.0₃0746'
$ws.Range("E38").Value = '  +4.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0400'
$ws.Range("E39").Value = '  +3.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '429.64'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.83'
$ws.Range("E41").Value = '  +4.86%  '
$ws.Range("D42").Value = '3.032.69'
$ws.Range("E42").Value = '  +6.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.51'
$ws.Range("E43").Value = '  +3.62%  '
$ws.Range("E44").Value = '  -4.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.274'
$ws.Range("E45").Value = '  +7.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.27'
$ws.Range("E46").Value = '  +8.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.49'
$ws.Range("E47").Value = '  +4.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.37'
$ws.Range("E48").Value = '  +4.42%  '
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '123.41'
$ws.Range("E51").Value = '  +4.27%  '
